$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.571607112884521
$ws.Range("B1").Value = 2.238519191741943
$ws.Range("C1").Value = 4.569770336151123
$ws.Range("D1").Value = 1.616632103919983
$ws.Range("E1").Value = 0.8047512769699097
